$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: replace the userName/password pair with a new set of test data
# (adds two new shared strings: "jakay34@gmail.com" and "12345678")
$ws.Range("A8").Value = "jakay34@gmail.com"

# Force B8 to be stored as text (quote-prefixed), matching the new
# "12345678" value being a text string rather than a number
$ws.Range("B8").Value = "'12345678"

# Update the active selection to D5
$ws.Range("D5").Select() | Out-Null
